# Add two new collected rows to the first sheet
# ("八位序列号收集收集结果yd5"), appending to the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("八位序列号收集收集结果yd5")

$ws.Range("A42").Value = "Satellite."
$ws.Range("B42").Value = 45915.9581712963
$ws.Range("C42").Value = "e272cd49"
$ws.Range("D42").Value = "'1125931910"

$ws.Range("A43").Value = "pots"
$ws.Range("B43").Value = 45916.0037152778
$ws.Range("C43").Value = "8c09b4b0"
$ws.Range("D43").Value = "'2014017597"

$ws.Range("B42:B43").NumberFormat = "yyyy/m/d h:mm:ss;@"

# The two serial-number cells above look like plain numbers, so Excel
# marks them with a "number stored as text" quote-prefix style when the
# literal value is entered. Re-apply the plain (un-prefixed) cell format
# from a sibling text cell so the stored style matches a normal text
# cell again, while keeping the text value itself intact.
$ws.Range("C42").Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("C43").Copy()
$ws.Range("D43").PasteSpecial(-4122)
$excel.CutCopyMode = $false
